$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50; this shifts all existing rows
# 50..136 down to 51..137, preserving their contents/styles.
$ws.Rows(50).Insert()

# Populate the newly inserted row 50 with the new weekly price record.
$ws.Range("A50").Value = 4
$ws.Range("B50").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C50").Value = "Los Lagos"
$ws.Range("D50").Value = 45028
$ws.Range("E50").Value = 10
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100104
$ws.Range("H50").Value = "Frutos de pepita"
$ws.Range("I50").Value = 100104003
$ws.Range("J50").Value = "Membrillo"
$ws.Range("K50").Value = "Champion"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 200
$ws.Range("N50").Value = 15000
$ws.Range("O50").Value = 16000
$ws.Range("P50").Value = 15500
$ws.Range("Q50").Value = '$/caja 18 kilos empedrada'
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 861
$ws.Range("T50").Value = 18
